# ============================================================
# Add 2022-Q4 data:
#  - 'sheet2' (currently named 2022-Q3) is duplicated; the copy
#    keeps the archived 2022-Q3 fund table and is renamed back
#    to '2022-Q3'. The original sheet is renamed to '2022-Q4'
#    and repopulated with the new quarter's fund table.
#  - The '总计' summary sheet gets a new row for 2022-Q4, and
#    the old 2022-Q3 summary row shifts down one row.
# ============================================================

$wb = $excel.ActiveWorkbook
$wsTotal = $wb.Worksheets.Item(1)   # "总计"
$wsQ3 = $wb.Worksheets.Item(2)      # currently "2022-Q3"

# ---- 1. Duplicate the current quarter sheet so its data is preserved ----
# as the new archived '2022-Q3' tab, placed right after it.
$wsQ3.Copy($null, $wsQ3)
$wsArchive = $wb.Worksheets.Item($wsQ3.Index + 1)

# ---- 2. Rename the original (now stale) sheet to the new quarter FIRST ----
# (frees up the "2022-Q3" name before the copy claims it)
$wsQ3.Name = "2022-Q4"
$wsQ4 = $wsQ3
$wsArchive.Name = "2022-Q3"

# ---- 3. Clear it out and rebuild with the 2022-Q4 fund table ----
$wsQ4.Cells.Clear()

# Match formatting used elsewhere in the workbook: header style + index style
# come from the '总计' sheet so style indices line up with the rest of the file.
$wsTotal.Range("B1:D1").Copy($wsQ4.Range("B1:H1"))
$wsTotal.Range("A2").Copy($wsQ4.Range("A2:A6"))

# Page margins matching the rest of the workbook
$wsQ4.PageSetup.LeftMargin = 54
$wsQ4.PageSetup.RightMargin = 54
$wsQ4.PageSetup.TopMargin = 72
$wsQ4.PageSetup.BottomMargin = 72
$wsQ4.PageSetup.HeaderMargin = 36
$wsQ4.PageSetup.FooterMargin = 36

# Header row
$wsQ4.Range("B1").Value = "基金代码"
$wsQ4.Range("C1").Value = "基金名称"
$wsQ4.Range("D1").Value = "基金规模"
$wsQ4.Range("E1").Value = "股票总仓位"
$wsQ4.Range("F1").Value = "仓位占比"
$wsQ4.Range("G1").Value = "持有市值(亿元)"
$wsQ4.Range("H1").Value = "仓位排名"

# Helper: force a numeric-looking string to be stored as TEXT (not auto-converted
# to a number), matching the source data where fund codes / percentages / sizes
# are text cells. The leading apostrophe forces text entry; resetting the style
# back to Normal afterwards drops the quote-prefix style iron_native applies.
function Set-TextValue($range, $value) {
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

# Row 2: 011001
$wsQ4.Range("A2").Value = 0
Set-TextValue $wsQ4.Range("B2") "011001"
$wsQ4.Range("C2").Value = "中邮兴荣价值一年持有期混合"
Set-TextValue $wsQ4.Range("D2") "5.44"
Set-TextValue $wsQ4.Range("E2") "67.79"
Set-TextValue $wsQ4.Range("F2") "2.21"
Set-TextValue $wsQ4.Range("G2") "0.1202"
$wsQ4.Range("H2").Value = 10

# Row 3: 010114
$wsQ4.Range("A3").Value = 1
Set-TextValue $wsQ4.Range("B3") "010114"
$wsQ4.Range("C3").Value = "华宝新兴成长混合A"
Set-TextValue $wsQ4.Range("D3") "3.09"
Set-TextValue $wsQ4.Range("E3") "89.74"
Set-TextValue $wsQ4.Range("F3") "2.77"
Set-TextValue $wsQ4.Range("G3") "0.0856"
$wsQ4.Range("H3").Value = 9

# Row 4: 240017
$wsQ4.Range("A4").Value = 2
Set-TextValue $wsQ4.Range("B4") "240017"
$wsQ4.Range("C4").Value = "华宝新兴产业混合"
Set-TextValue $wsQ4.Range("D4") "2.74"
Set-TextValue $wsQ4.Range("E4") "86.64"
Set-TextValue $wsQ4.Range("F4") "2.72"
Set-TextValue $wsQ4.Range("G4") "0.0745"
$wsQ4.Range("H4").Value = 9

# Row 5: 007305
$wsQ4.Range("A5").Value = 3
Set-TextValue $wsQ4.Range("B5") "007305"
$wsQ4.Range("C5").Value = "国联安新科技混合"
Set-TextValue $wsQ4.Range("D5") "1.47"
Set-TextValue $wsQ4.Range("E5") "77.63"
Set-TextValue $wsQ4.Range("F5") "2.48"
Set-TextValue $wsQ4.Range("G5") "0.0365"
$wsQ4.Range("H5").Value = 9

# Row 6: 017197
$wsQ4.Range("A6").Value = 4
Set-TextValue $wsQ4.Range("B6") "017197"
$wsQ4.Range("C6").Value = "华宝新兴成长混合C"
Set-TextValue $wsQ4.Range("D6") "0.00"
Set-TextValue $wsQ4.Range("E6") "89.74"
Set-TextValue $wsQ4.Range("F6") "2.77"
$wsQ4.Range("G6").Value = 0
$wsQ4.Range("H6").Value = 9

# ---- 4. '总计' summary sheet: insert the 2022-Q4 row, push 2022-Q3 down ----

# Copy A2's style down to A3 first so the shifted row keeps its formatting
$wsTotal.Range("A2").Copy($wsTotal.Range("A3"))

# Shift the existing 2022-Q3 summary row down to row 3
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q3"
$wsTotal.Range("C3").Value = 1
$wsTotal.Range("D3").Value = 0.14

# Write the new 2022-Q4 summary row
$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 5
$wsTotal.Range("D2").Value = 0.32

# Leave the workbook focused on the summary sheet (matches the book-level
# activeTab="0" the workbook already had before this edit).
$wsTotal.Activate()

